$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 110, shifting existing rows 110-115 down to 111-116
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new weekly record
$ws.Cells.Item(110, 1).Value = 1
$ws.Cells.Item(110, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(110, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(110, 4).Value = 44783
$ws.Cells.Item(110, 5).Value = 15
$ws.Cells.Item(110, 6).Value = 100112008
$ws.Cells.Item(110, 7).Value = "Coliflor"
$ws.Cells.Item(110, 8).Value = "Sin especificar"
$ws.Cells.Item(110, 9).Value = "Tercera"
$ws.Cells.Item(110, 10).Value = 900
$ws.Cells.Item(110, 11).Value = 500
$ws.Cells.Item(110, 12).Value = 600
$ws.Cells.Item(110, 13).Value = 550
$ws.Cells.Item(110, 14).Value = "$/unidad"
$ws.Cells.Item(110, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(110, 16).Value = 550
$ws.Cells.Item(110, 17).Value = 1
$ws.Cells.Item(110, 18).Value = "Hortaliza"
